$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the capitalization of the "etoposido" drug name to "Etoposido"
# (this also merges with the existing duplicate "Etoposido" shared string)
$ws.Range("E4").Value = "Etoposido"

# Leave the selection on the cell that was last interacted with
$ws.Range("B9").Activate()
